# "almost completed corporate master"
#
# A new first data field - "Equity Type" - is introduced into the Sample
# Folios header row. In the source data this shows up as the header row
# (row 1) shifting one column to the right (Folio -> B, Certificate
# Number -> C, ... , Share Certificate Master Id* -> N) with the brand
# new "Equity Type" header taking over column A, and the matching
# second row (the "template"/format row under the headers) carrying its
# per-column formatting one column to the right as well, so each header
# keeps the cell format it originally had. The remaining body rows
# (3-10) simply gain a 14th column that copies the formatting already
# used by column M in that row.
#
# NOTE: `.Value` getters are unreliable on this host - use `.Value2` for
# reads; `.Value2` works fine as a setter too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$lastCol = 14  # N - one more than the old M (13)

# 1. Shift the header row (row 1) one column to the right, working from
#    the rightmost column back towards B so we never clobber a value we
#    still need to read.
for ($col = $lastCol; $col -ge 2; $col--) {
    $ws.Cells.Item(1, $col).Value2 = $ws.Cells.Item(1, $col - 1).Value2
}
$ws.Cells.Item(1, 1).Value2 = "Equity Type"

# Column N's header cell is brand new - give it the same format column M's
# header used to have (they're all the same style, but be explicit).
$ws.Cells.Item(1, 13).Copy() | Out-Null
$ws.Cells.Item(1, 14).PasteSpecial($xlPasteFormats) | Out-Null

# 2. Shift row 2's per-column formatting one column to the right the same
#    way, so "Equity Type" (now column A) ends up formatted like "Folio"
#    used to be, "Folio" (now column B) keeps its own old look, etc.
for ($col = $lastCol; $col -ge 2; $col--) {
    $ws.Cells.Item(2, $col - 1).Copy() | Out-Null
    $ws.Cells.Item(2, $col).PasteSpecial($xlPasteFormats) | Out-Null
}

# 3. Rows 3-10 just need a new column N that matches column M's format.
for ($row = 3; $row -le 10; $row++) {
    $ws.Cells.Item($row, 13).Copy() | Out-Null
    $ws.Cells.Item($row, 14).PasteSpecial($xlPasteFormats) | Out-Null
}
